$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix the "Deplacement des pages lisp, pascal, PHP, types_langage," entry:
#        drop the stray trailing comma (row 54, column B / Description).
$ws.Cells.Item(54, 2).Value = "Déplacement des pages lisp, pascal, PHP, types_langage"

# --- 2. Replace the placeholder "?" software entry for the robotique page creation
#        (row 59, column E / Logiciel) with the real value.
$ws.Cells.Item(59, 5).Value = "Brackets"

# --- 3. Append the remaining journal entries (rows 79-86), copying the formatting
#        from the last existing row (78) so styles/number formats line up.
$ws.Range("A78:F78").Copy() | Out-Null
$ws.Range("A79:F86").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

$rows = @(
    @{ Row = 79; A = 44202; B = "Ajout d'une image JPG shéma monde IA"; C = "--------"; D = "Elisa"; E = "Github Web"; F = "--------" },
    @{ Row = 80; A = 44202; B = "Modification de la feuille de style générale"; C = "Changement de la couleur de la classe code"; D = "Mathieu"; E = "VS Code"; F = "CSS" },
    @{ Row = 81; A = 44202; B = "Modification de la page programmation"; C = "Remise en forme du code"; D = "Mathieu"; E = "VS Code"; F = "HTML" },
    @{ Row = 82; A = 44202; B = "Renomination de l'image JPG shéma monde IA"; C = "Non respect des conventions (espaces)"; D = "Mathieu"; E = "--------"; F = "--------" },
    @{ Row = 83; A = 44202; B = "Modification du script JS"; C = "Remise en forme du code et ajout d'instructions pour la fonction info"; D = "Mathieu"; E = "VS Code"; F = "Javascript" },
    @{ Row = 84; A = 43841; B = "Modification de la feuille de style générale"; C = "Transformation des styles page a en a"; D = "Mathieu"; E = "VS Code"; F = "CSS" },
    @{ Row = 85; A = 43841; B = "Modification du script JS"; C = "Création de la fonction explication"; D = "Mathieu"; E = "VS Code"; F = "Javascript" },
    @{ Row = 86; A = 43841; B = "Modification de la page programmation"; C = "Insertion d'appel à la fct explication pr certains mots"; D = "Mathieu"; E = "VS Code"; F = "HTML" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
}

# Rows 83 and 86 wrap onto two lines in the source workbook (taller row height).
$ws.Rows.Item(83).RowHeight = 30
$ws.Rows.Item(86).RowHeight = 30

# --- 4. Re-apply the existing header freeze (row 1) so the frozen pane survives
#        the new rows, then leave the final selection on the newly added last
#        row, matching the author's save state.
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $false
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("D86:F86").Select() | Out-Null
